$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.017795677743293
$ws.Range("D2").Value = 1.054144836612939
$ws.Range("E2").Value = 1.019144894807285
$ws.Range("F2").Value = 1.056101277921388
$ws.Range("I2").Value = 1.043581809496192
$ws.Range("J2").Value = 1.023007760634219
$ws.Range("K2").Value = 1.056888855380177
$ws.Range("L2").Value = 1.021989012044209
$ws.Range("M2").Value = 1.058839913578404
$ws.Range("N2").Value = 1.02446054899495

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.018915747269516
$ws.Range("D3").Value = 1.054879735223907
$ws.Range("E3").Value = 1.020098387882625
$ws.Range("F3").Value = 1.057047306033185
$ws.Range("I3").Value = 1.043784867488091
$ws.Range("J3").Value = 1.023762977472822
$ws.Range("K3").Value = 1.057437139279821
$ws.Range("L3").Value = 1.022748061261864
$ws.Range("M3").Value = 1.059599176711343
$ws.Range("N3").Value = 1.025216838328088

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.019639890275314
$ws.Range("D4").Value = 1.055352038979693
$ws.Range("E4").Value = 1.020715235621115
$ws.Range("F4").Value = 1.057655893421346
$ws.Range("I4").Value = 1.043913004342863
$ws.Range("J4").Value = 1.024250541801667
$ws.Range("K4").Value = 1.057788029705283
$ws.Range("L4").Value = 1.023238478553621
$ws.Range("M4").Value = 1.060086298588383
$ws.Range("N4").Value = 1.025705095054204

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.019944173496857
$ws.Range("D5").Value = 1.055549822209866
$ws.Range("E5").Value = 1.020974528493887
$ws.Range("F5").Value = 1.05791089141766
$ws.Range("I5").Value = 1.043966092931432
$ws.Range("J5").Value = 1.024455248107387
$ws.Range("K5").Value = 1.057934612540079
$ws.Range("L5").Value = 1.023444473516465
$ws.Range("M5").Value = 1.06029008477708
$ws.Range("N5").Value = 1.025910092066358

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.019995255431676
$ws.Range("D6").Value = 1.055582985501679
$ws.Range("E6").Value = 1.021018063160744
$ws.Range("F6").Value = 1.057953656707623
$ws.Range("I6").Value = 1.04397496096574
$ws.Range("J6").Value = 1.02448960362332
$ws.Range("K6").Value = 1.057959169782794
$ws.Range("L6").Value = 1.023479050636778
$ws.Range("M6").Value = 1.060324242711335
$ws.Range("N6").Value = 1.025944496371065

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.019643956696854
$ws.Range("D7").Value = 1.055354684809116
$ws.Range("E7").Value = 1.020718700422477
$ws.Range("F7").Value = 1.057659304066697
$ws.Range("I7").Value = 1.043913716781806
$ws.Range("J7").Value = 1.024253278141165
$ws.Range("K7").Value = 1.0577899920142
$ws.Range("L7").Value = 1.02324123176151
$ws.Range("M7").Value = 1.060089025519713
$ws.Range("N7").Value = 1.025707835279619

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.018174338934345
$ws.Range("D8").Value = 1.054393866027924
$ws.Range("E8").Value = 1.01946715851846
$ws.Range("F8").Value = 1.056421728398157
$ws.Range("I8").Value = 1.043651107429434
$ws.Range("J8").Value = 1.023263220058015
$ws.Range("K8").Value = 1.057074954398719
$ws.Range("L8").Value = 1.022245689334404
$ws.Range("M8").Value = 1.059097373390009
$ws.Range("N8").Value = 1.02471637120043

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.015579904988492
$ws.Range("D9").Value = 1.05267613290374
$ws.Range("E9").Value = 1.017260802782369
$ws.Range("F9").Value = 1.054213756800555
$ws.Range("I9").Value = 1.043163458913932
$ws.Range("J9").Value = 1.021510080796128
$ws.Range("K9").Value = 1.055785253697704
$ws.Range("L9").Value = 1.020485746968578
$ws.Range("M9").Value = 1.057318041700642
$ws.Range("N9").Value = 1.022960742279703

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.013846986316556
$ws.Range("D10").Value = 1.051514464936902
$ws.Range("E10").Value = 1.01578921638154
$ws.Range("F10").Value = 1.05272351364828
$ws.Range("I10").Value = 1.042821661333204
$ws.Range("J10").Value = 1.02033555565162
$ws.Range("K10").Value = 1.054905532800946
$ws.Range("L10").Value = 1.019308614710917
$ws.Range("M10").Value = 1.05611040622087
$ws.Range("N10").Value = 1.02178454917477

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.013095811690957
$ws.Range("D11").Value = 1.051007543170722
$ws.Range("E11").Value = 1.015151834351867
$ws.Range("F11").Value = 1.052073893087153
$ws.Range("I11").Value = 1.04266970568987
$ws.Range("J11").Value = 1.019825597486589
$ws.Range("K11").Value = 1.054519888572751
$ws.Range("L11").Value = 1.01879798616977
$ws.Range("M11").Value = 1.055582413307307
$ws.Range("N11").Value = 1.02127386681065

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01281666850997
$ws.Range("D12").Value = 1.050818662677487
$ws.Range("E12").Value = 1.014915055186171
$ws.Range("F12").Value = 1.051831943613893
$ws.Range("I12").Value = 1.042612668818649
$ws.Range("J12").Value = 1.019635967969063
$ws.Range("K12").Value = 1.054375934630188
$ws.Range("L12").Value = 1.018608176646249
$ws.Range("M12").Value = 1.055385530344435
$ws.Range("N12").Value = 1.021083967997459

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.01287655128804
$ws.Range("D13").Value = 1.050859204747765
$ws.Range("E13").Value = 1.014965846359671
$ws.Range("F13").Value = 1.051883872080686
$ws.Range("I13").Value = 1.042624930284666
$ws.Range("J13").Value = 1.01967665358679
$ws.Range("K13").Value = 1.054406845293494
$ws.Range("L13").Value = 1.018648897751535
$ws.Range("M13").Value = 1.055427796929195
$ws.Range("N13").Value = 1.02112471139343

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.01307274015845
$ws.Range("D14").Value = 1.050991942225698
$ws.Range("E14").Value = 1.01513226266945
$ws.Range("F14").Value = 1.052053906759121
$ws.Range("I14").Value = 1.042665003114971
$ws.Range("J14").Value = 1.01980992690691
$ws.Range("K14").Value = 1.054508003745581
$ws.Range("L14").Value = 1.018782299302455
$ws.Range("M14").Value = 1.055566154462347
$ws.Range("N14").Value = 1.021258173976952

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.013193602177676
$ws.Range("D15").Value = 1.051073648374484
$ws.Range("E15").Value = 1.015234793636345
$ws.Range("F15").Value = 1.052158584404591
$ws.Range("I15").Value = 1.042689614638476
$ws.Range("J15").Value = 1.019892013363221
$ws.Range("K15").Value = 1.054570236936647
$ws.Range("L15").Value = 1.01886447391285
$ws.Range("M15").Value = 1.055651300017714
$ws.Range("N15").Value = 1.021340377005446

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.013896822090082
$ws.Range("D16").Value = 1.051548025294007
$ws.Range("E16").Value = 1.015831513568996
$ws.Range("F16").Value = 1.052766535474562
$ws.Range("I16").Value = 1.042831662837692
$ws.Range("J16").Value = 1.020369370725824
$ws.Range("K16").Value = 1.054931027354156
$ws.Range("L16").Value = 1.019342483945915
$ws.Range("M16").Value = 1.056145340382042
$ws.Range("N16").Value = 1.021818412270257

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.014337715350257
$ws.Range("D17").Value = 1.051844542488751
$ws.Range("E17").Value = 1.016205772739208
$ws.Range("F17").Value = 1.053146726344366
$ws.Range("I17").Value = 1.042919707356964
$ws.Range("J17").Value = 1.020668434048538
$ws.Range("K17").Value = 1.055156078555013
$ws.Range("L17").Value = 1.019642079509708
$ws.Range("M17").Value = 1.056453879249306
$ws.Range("N17").Value = 1.022117900297193

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.014594802961588
$ws.Range("D18").Value = 1.052017118481605
$ws.Range("E18").Value = 1.016424054908972
$ws.Range("F18").Value = 1.053368066952712
$ws.Range("I18").Value = 1.04297068060641
$ws.Range("J18").Value = 1.020842739435941
$ws.Range("K18").Value = 1.055286891741444
$ws.Range("L18").Value = 1.019816739577082
$ws.Range("M18").Value = 1.056633354717508
$ws.Range("N18").Value = 1.022292453218238

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.014682450025065
$ws.Range("D19").Value = 1.052075898392403
$ws.Range("E19").Value = 1.01649848064836
$ws.Range("F19").Value = 1.053443467432576
$ws.Range("I19").Value = 1.042987996422813
$ws.Range("J19").Value = 1.020902150503215
$ws.Range("K19").Value = 1.055331418401831
$ws.Range("L19").Value = 1.019876279080302
$ws.Range("M19").Value = 1.056694468090196
$ws.Range("N19").Value = 1.022351948656043

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.014290419753731
$ws.Range("D20").Value = 1.051812768039093
$ws.Range("E20").Value = 1.016165620022464
$ws.Range("F20").Value = 1.053105978728418
$ws.Range("I20").Value = 1.04291030048365
$ws.Range("J20").Value = 1.020636361166604
$ws.Range("K20").Value = 1.055131979788034
$ws.Range("L20").Value = 1.019609944947351
$ws.Range("M20").Value = 1.056420826633572
$ws.Range("N20").Value = 1.022085781868088

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.013014970849856
$ws.Range("D21").Value = 1.050952870551943
$ws.Range("E21").Value = 1.015083257968304
$ws.Range("F21").Value = 1.05200385376961
$ws.Range("I21").Value = 1.042653219053318
$ws.Range("J21").Value = 1.019770686994467
$ws.Range("K21").Value = 1.054478234657245
$ws.Range("L21").Value = 1.018743019723521
$ws.Range("M21").Value = 1.055525432676012
$ws.Range("N21").Value = 1.021218878339332

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.012212330135573
$ws.Range("D22").Value = 1.050408821844668
$ws.Range("E22").Value = 1.01440257679786
$ws.Range("F22").Value = 1.051307134715074
$ws.Range("I22").Value = 1.042488146503015
$ws.Range("J22").Value = 1.019225197385895
$ws.Range("K22").Value = 1.054063100148581
$ws.Range("L22").Value = 1.0181971428102
$ws.Range("M22").Value = 1.054958049135387
$ws.Range("N22").Value = 1.020672614072947

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.012637893609061
$ws.Range("D23").Value = 1.050697554199894
$ws.Range("E23").Value = 1.014763433923098
$ws.Range("F23").Value = 1.051676836025604
$ws.Range("I23").Value = 1.042575980079856
$ws.Range("J23").Value = 1.019514486352124
$ws.Range("K23").Value = 1.054283559194046
$ws.Range("L23").Value = 1.018486599207593
$ws.Range("M23").Value = 1.055259248242632
$ws.Range("N23").Value = 1.02096231386269

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.014311790838282
$ws.Range("D24").Value = 1.051827126711154
$ws.Range("E24").Value = 1.01618376335858
$ws.Range("F24").Value = 1.053124392114003
$ws.Range("I24").Value = 1.042914552223774
$ws.Range("J24").Value = 1.020650853932531
$ws.Range("K24").Value = 1.055142870390197
$ws.Range("L24").Value = 1.019624465447375
$ws.Range("M24").Value = 1.056435763201508
$ws.Range("N24").Value = 1.022100295215404

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.016251203118025
$ws.Range("D25").Value = 1.053123122546068
$ws.Range("E25").Value = 1.017831317023763
$ws.Range("F25").Value = 1.054787790942505
$ws.Range("I25").Value = 1.043292473788352
$ws.Range("J25").Value = 1.021964323192706
$ws.Range("K25").Value = 1.056122187205804
$ws.Range("L25").Value = 1.020941408618346
$ws.Range("M25").Value = 1.057781818763181
$ws.Range("N25").Value = 1.023415629752586
